# Price tracker update: insert a new snapshot column before column B,
# shifting all existing date columns (B:Q) one column to the right
# (C:R) and filling the new column B with the latest price snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column B; this shifts B:Q -> C:R
# and copies formatting from column B (so the header style carries over).
$ws.Columns("B").Insert()

# Match the width used by all the other data columns (the OOXML "width"
# unit of 21 corresponds to a ColumnWidth property value of ~20.14 given
# this engine's character-width <-> pixel rounding).
$ws.Columns("B").ColumnWidth = 20.14

# New snapshot timestamp for the header row.
$ws.Cells.Item(1, 2).Value = "2025-12-22 00:24"

# New price snapshot values for column B (row 2 .. row 26).
# $null / empty entries mean the cell stays blank (no price recorded).
$prices = @(929, $null, $null, 569, 499, $null, $null, $null, $null, 2997, 569, $null, $null, 499, $null, 929, 499, 1497, 929, 499, $null, 1299, $null, $null, $null)

for ($i = 0; $i -lt $prices.Length; $i++) {
    $row = $i + 2
    $val = $prices[$i]
    if ($null -ne $val) {
        $ws.Cells.Item($row, 2).Value = $val
    }
}
